# Add a new reference row (row 79) to the data-reference sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

# Columns A and C hold numeric-looking codes ("015", "112233") that must be
# stored as text (so a leading zero like "015" is preserved). Pre-formatting
# the cells as Text before assigning the value keeps Excel from
# auto-converting the string into a number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "015"

$ws.Cells.Item($row, 2).Value = "Kementerian Keuangan"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "112233"

$ws.Cells.Item($row, 4).Value = "TESTT"
$ws.Cells.Item($row, 5).Value = "TESTT TESTTTTTTTTTT"

# Column F is a plain running-count number.
$ws.Cells.Item($row, 6).Value = 78
